$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "tvIigywR"
$ws.Range("B3").Value = "31/10/2024"
$ws.Range("C3").Value = "11:30"
$ws.Range("D3").Value = "QATAR - QSL"
$ws.Range("E3").Value = "Al Khor"
$ws.Range("F3").Value = "Qatar SC"
$ws.Range("G3").Value = 2.72
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 2.45
$ws.Range("J3").Value = 3.3
$ws.Range("K3").Value = 2.12
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 7.8
$ws.Range("O3").Value = 1.24
$ws.Range("P3").Value = 3.65
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 2.02
$ws.Range("S3").Value = 1.37
$ws.Range("T3").Value = 2.85
$ws.Range("U3").Value = 1.57
$ws.Range("V3").Value = 2.25
$ws.Range("W3").Value = 10.25
$ws.Range("X3").Value = 15.5
$ws.Range("Y3").Value = 9.75
$ws.Range("Z3").Value = 32
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 7.8
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 11.5
$ws.Range("AF3").Value = 45
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 14
$ws.Range("AJ3").Value = 9.25
$ws.Range("AK3").Value = 27
$ws.Range("AL3").Value = 18.5
$ws.Range("AM3").Value = 23
$ws.Range("AN3").Value = 4.85
$ws.Range("AO3").Value = 15
$ws.Range("AP3").Value = 20
$ws.Range("AQ3").Value = 65
$ws.Range("AR3").Value = 90
$ws.Range("AS3").Value = 250
$ws.Range("AT3").Value = 2.85
$ws.Range("AU3").Value = 6.5
$ws.Range("AV3").Value = 50
$ws.Range("AW3").Value = 4.55
$ws.Range("AX3").Value = 13
$ws.Range("AY3").Value = 18.5
$ws.Range("AZ3").Value = 50
$ws.Range("BA3").Value = 75
$ws.Range("BB3").Value = 200
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51
$ws.Range("A4").Value = "OxptyxHr"
$ws.Range("B4").Value = "31/10/2024"
$ws.Range("C4").Value = "11:30"
$ws.Range("D4").Value = "QATAR - QSL"
$ws.Range("E4").Value = "Al-Gharafa"
$ws.Range("F4").Value = "Al Arabi"
$ws.Range("G4").Value = 2.12
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 2.45
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 9.75
$ws.Range("O4").Value = 1.13
$ws.Range("P4").Value = 5.2
$ws.Range("Q4").Value = 1.4
$ws.Range("R4").Value = 2.7
$ws.Range("S4").Value = 1.23
$ws.Range("T4").Value = 3.7
$ws.Range("U4").Value = 1.38
$ws.Range("V4").Value = 2.82
$ws.Range("W4").Value = 13.5
$ws.Range("X4").Value = 14.5
$ws.Range("Y4").Value = 9.25
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 14.5
$ws.Range("AB4").Value = 17
$ws.Range("AC4").Value = 9.75
$ws.Range("AD4").Value = 8.5
$ws.Range("AE4").Value = 10.5
$ws.Range("AF4").Value = 28
$ws.Range("AG4").Value = 120
$ws.Range("AH4").Value = 17
$ws.Range("AI4").Value = 22
$ws.Range("AJ4").Value = 11
$ws.Range("AK4").Value = 40
$ws.Range("AL4").Value = 20
$ws.Range("AM4").Value = 19.5
$ws.Range("AN4").Value = 4.7
$ws.Range("AO4").Value = 10.5
$ws.Range("AP4").Value = 13.5
$ws.Range("AQ4").Value = 35
$ws.Range("AR4").Value = 45
$ws.Range("AS4").Value = 110
$ws.Range("AT4").Value = 3.7
$ws.Range("AU4").Value = 6
$ws.Range("AV4").Value = 32
$ws.Range("AW4").Value = 5.6
$ws.Range("AX4").Value = 14.5
$ws.Range("AY4").Value = 15.5
$ws.Range("AZ4").Value = 55
$ws.Range("BA4").Value = 55
$ws.Range("BB4").Value = 120
$ws.Range("BC4").Value = 350
$ws.Range("BD4").Value = 51
